$wb = $excel.ActiveWorkbook
foreach ($ws in $wb.Worksheets) {
    Write-Output $ws.Name
}
